$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Day 16" / "Day 20" columns already present further right in the sheet
# (O/R = first occurrence, AN/AQ = second occurrence) get copied into the
# front B:E columns, replacing the old 1/2/3/4 placeholder header values and
# the old sample data that lived there. Where the source column had no value
# for a subject, the destination cell is cleared instead of overwritten.

# Row 1 (headers)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON)
$ws.Range("B2").Value = 20.804287499999987
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 29.091802499999972
$ws.Range("E2").Value = 29.008192500000007

# Row 3 (STR)
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 31.909106249999979
$ws.Range("D3").Value = 25.356307499999957
$ws.Range("E3").Value = 28.058771999999976

# Selection now only spans the still-populated sample columns.
$ws.Range("B1:E3").Select() | Out-Null
